# Staging.LocationType.xlsx header-row relabel.
#
# The sheet's header row (row 2) lists the staging columns. The shared
# strings backing it were reshuffled upstream so that, with the cell
# references left untouched, the labels that now line up with columns
# A-D are: Code, Description, LocationType_ID, Name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "Description"
$ws.Range("C2").Value = "LocationType_ID"
$ws.Range("D2").Value = "Name"
